$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds text-formatted numbers (e.g. thousands separated
# by dots like "35.371.39", or values with significant trailing zeros like
# "7.90"). Force these cells to Text format first so Excel does not
# reinterpret/renormalize them as numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.371.39"
$ws.Range("E2").Value = "  +2.69%  "
$ws.Range("D3").Value = "1.844.94"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("D5").Value = "230.42"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "42.74"
$ws.Range("E8").Value = "  +12.12%  "
$ws.Range("E9").Value = "  +6.96%  "
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D12").Value = "2.112.93"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("D13").Value = "1.847.02"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("D14").Value = "11.31"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").Value = "0.671"
$ws.Range("E15").Value = "  +7.21%  "
$ws.Range("E16").Value = "  +6.59%  "
$ws.Range("D17").Value = "35.390.29"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "70.27"
$ws.Range("E18").Value = "  +3.54%  "
$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").Value = "244.97"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "12.09"
$ws.Range("E21").Value = "  +9.59%  "
$ws.Range("D22").Value = "4.63"
$ws.Range("E22").Value = "  +13.38%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").Value = "168.99"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "7.90"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("D27").Value = "17.75"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("E29").Value = "  +13.69%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "3.341.47"
$ws.Range("E31").Value = "  +37.53%  "
$ws.Range("E32").Value = "  +6.34%  "
$ws.Range("D33").Value = "3.93"
$ws.Range("E33").Value = "  +4.99%  "
$ws.Range("E34").Value = "  +5.51%  "
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").Value = "96.22"
$ws.Range("E36").Value = "  +16.80%  "
$ws.Range("E37").Value = "  +7.75%  "
$ws.Range("D38").Value = "1.345.37"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").Value = "1.09"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("D41").Value = "0.0194"
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("E42").Value = "  +6.29%  "
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("D44").Value = "14.76"
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "6.24"
$ws.Range("E47").Value = "  +8.81%  "
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").Value = "2.013.00"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "103.43"
$ws.Range("E51").Value = "  +1.40%  "
